$wb = $excel.ActiveWorkbook

# Sheet R1
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3927:21:57"
$ws.Range("G3").Value = "66:54:35"

# Sheet R2
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12108:45:38"
$ws.Range("G3").Value = "3238:29:07"
$ws.Range("G4").Value = "476:40:41"

# Sheet R4
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2954:35:27"
$ws.Range("G3").Value = "181:47:42"

# Sheet R5
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "428:34:26"

# Sheet R6
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "69:06:44"
